$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a number need to be pinned as
# text (NumberFormat "@") while writing, then have formatting
# cleared again afterward so no stray style survives the write -
# matching the original "General"-styled inline-string cells.

$ws.Range('D2').Value = '28.320.96'
$ws.Range('E2').Value = '  +3.70%  '
$ws.Range('D3').Value = '1.786.37'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3844'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3443'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.06'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.154'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07416'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.54'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +8.20%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.465'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.378'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.26%  '
$ws.Range('D16').Value = '1.785.33'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001080'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06702'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.51'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.429'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '28.307.53'
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.10'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.367'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.79'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.416'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '154.59'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.76'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('B31').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C31').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D31').Value = '1.985.64'
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.128'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.006'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08971'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.75'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02424'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6885'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.368'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06410'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2158'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.250'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.503'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -6.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.305'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.19'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9999'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6299'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.881'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.90'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.084'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07482'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.01%  '
$ws.Range('E51').Value = '  +6.42%  '
